# Rename the inline picture shapes' docPr/name attributes:
#   - Footer (primary)   logo: image2.png -> image1.png  (wp:docPr id="1")
#   - Footer (first page) logo: image2.png -> image1.png (wp:docPr id="2")
#   - Header (first page) logo: image1.jpg -> image2.jpg (wp:docPr id="3")
#
# InlineShape has no settable .Name in the Word object model (only the
# floating Shape object exposes it), so each picture is round-tripped
# through ConvertToShape()/ConvertToInlineShape() to rename it while
# keeping it inline in the header/footer flow.

$d = $word.ActiveDocument

function Rename-InlinePicture($inlineShape, $newName) {
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    [void]$floating.ConvertToInlineShape()
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    Rename-InlinePicture $shp "image1.png"
                }
            }
        }
    }

    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-InlinePicture $shp "image2.jpg"
                }
            }
        }
    }
}
